$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate/swap existing rows to their new team/odds data ---
function Set-RowFV($row, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 6 + $i).Value = $vals[$i]
    }
}

Set-RowFV 73 @("Juventude", 3, "Atletico GO", 0, 2.5, "21/05/2023 16:12", 2.15, "24/05/2023 23:59", 3.12, "21/05/2023 16:12", 3.35, "24/05/2023 23:59", 3.1, "21/05/2023 16:12", 3.7, "24/05/2023 23:59", "https://www.betexplorer.com/football/brazil/serie-b/esporte-clube-juventude-atletico-go/ljr1zNpG/")
Set-RowFV 74 @("Vila Nova FC", 1, "Ituano", 0, 1.84, "21/05/2023 23:42", 1.77, "24/05/2023 23:56", 3.32, "21/05/2023 23:42", 3.39, "24/05/2023 23:53", 4.98, "21/05/2023 23:42", 5.65, "24/05/2023 23:56", "https://www.betexplorer.com/football/brazil/serie-b/vila-nova-fc-ituano/428cG1gp/")
Set-RowFV 76 @("Vitoria", 1, "CRB", 0, 1.78, "21/05/2023 23:42", 1.71, "24/05/2023 23:37", 3.35, "21/05/2023 23:42", 3.41, "24/05/2023 23:59", 4.92, "21/05/2023 23:42", 5.69, "24/05/2023 23:59", "https://www.betexplorer.com/football/brazil/serie-b/vitoria-crb/SzfhxqG3/")
Set-RowFV 307 @("Sampaio Correa", 1, "Novorizontino", 1, 3.32, "02/10/2023 06:12", 3.03, "07/10/2023 21:59", 2.94, "02/10/2023 06:12", 2.89, "07/10/2023 21:51", 2.39, "02/10/2023 06:12", 2.8, "07/10/2023 21:59", "https://www.betexplorer.com/football/brazil/serie-b/sampaio-correa-novorizontino/S4ibIC04/")
Set-RowFV 308 @("Botafogo SP", 0, "Avai", 1, 2.02, "02/10/2023 06:12", 2.05, "07/10/2023 21:58", 3.01, "02/10/2023 06:12", 3.04, "07/10/2023 21:58", 4.62, "02/10/2023 06:12", 4.58, "07/10/2023 21:58", "https://www.betexplorer.com/football/brazil/serie-b/botafogo-sp-avai/hnaoLEGo/")
Set-RowFV 334 @("Novorizontino", 2, "Ponte Preta", 0, 1.53, "24/10/2023 01:12", 1.39, "28/10/2023 02:21", 3.81, "24/10/2023 01:12", 4.57, "28/10/2023 02:28", 6.57, "24/10/2023 01:12", 10.13, "28/10/2023 02:28", "https://www.betexplorer.com/football/brazil/serie-b/novorizontino-ponte-preta/EeWldSr9/")
Set-RowFV 335 @("Guarani", 0, "Botafogo SP", 1, 1.64, "21/10/2023 22:12", 1.7, "28/10/2023 02:05", 3.41, "21/10/2023 22:12", 3.43, "28/10/2023 02:17", 6.14, "21/10/2023 22:12", 6.36, "28/10/2023 02:24", "https://www.betexplorer.com/football/brazil/serie-b/guarani-botafogo-sp/QZkfzoMk/")
Set-RowFV 336 @("Ituano", 0, "Mirassol", 0, 3.08, "24/10/2023 01:12", 3.71, "28/10/2023 20:24", 2.89, "24/10/2023 01:12", 2.9, "28/10/2023 20:24", 2.57, "24/10/2023 01:12", 2.39, "28/10/2023 20:24", "https://www.betexplorer.com/football/brazil/serie-b/ituano-mirassol/z5Xpc8T2/")
Set-RowFV 337 @("Criciuma", 3, "Sampaio Correa", 0, 1.7, "21/10/2023 23:42", 1.58, "28/10/2023 20:26", 3.42, "21/10/2023 23:42", 3.72, "28/10/2023 20:26", 5.95, "21/10/2023 23:42", 7.22, "28/10/2023 20:26", "https://www.betexplorer.com/football/brazil/serie-b/criciuma-sampaio-correa/bml2YOj2/")


# --- Append two new match rows (352, 353), copying formatting from row 351 ---
$ws.Range("A351:V351").Copy()
$ws.Range("A352:V353").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

function Set-RowAV($row, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 1 + $i
        $cell = $ws.Cells.Item($row, $col)
        if ($col -eq 4) {
            # column D ("temporada") is a text value that looks numeric (e.g. "2023");
            # force text so Excel doesn't silently coerce it to a number.
            $cell.NumberFormat = "@"
            $cell.Value = $vals[$i]
            $cell.ClearFormats()
        } else {
            $cell.Value = $vals[$i]
        }
    }
}

Set-RowAV 352 @(351, "brazil", "serie-b", "2023", 45240.95833333334, "Vila Nova FC", 4, "Londrina", 1, 1.54, "05/11/2023 22:12", 1.41, "10/11/2023 22:54", 3.78, "05/11/2023 22:12", 4.54, "10/11/2023 22:58", 7.43, "05/11/2023 22:12", 9.19, "10/11/2023 22:58", "https://www.betexplorer.com/football/brazil/serie-b/vila-nova-fc-londrina/ITIOoZQI/")
Set-RowAV 353 @(352, "brazil", "serie-b", "2023", 45241.0625, "Sport Recife", 0, "Atletico GO", 0, 1.81, "04/11/2023 06:43", 1.92, "11/11/2023 01:25", 3.51, "04/11/2023 06:43", 3.34, "11/11/2023 01:25", 4.82, "04/11/2023 06:43", 4.69, "11/11/2023 01:25", "https://www.betexplorer.com/football/brazil/serie-b/sport-recife-atletico-go/WzZvlP7n/")

